$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original row values (columns B:AB) for all affected rows before any writes
$row41 = $ws.Range("B41:AB41").Value2
$row42 = $ws.Range("B42:AB42").Value2
$row43 = $ws.Range("B43:AB43").Value2
$row53 = $ws.Range("B53:AB53").Value2
$row54 = $ws.Range("B54:AB54").Value2
$row80 = $ws.Range("B80:AB80").Value2
$row81 = $ws.Range("B81:AB81").Value2
$row93 = $ws.Range("B93:AB93").Value2
$row94 = $ws.Range("B94:AB94").Value2
$row98 = $ws.Range("B98:AB98").Value2
$row99 = $ws.Range("B99:AB99").Value2
$row102 = $ws.Range("B102:AB102").Value2
$row103 = $ws.Range("B103:AB103").Value2
$row111 = $ws.Range("B111:AB111").Value2
$row112 = $ws.Range("B112:AB112").Value2
$row123 = $ws.Range("B123:AB123").Value2
$row124 = $ws.Range("B124:AB124").Value2
$row130 = $ws.Range("B130:AB130").Value2
$row132 = $ws.Range("B132:AB132").Value2
$row137 = $ws.Range("B137:AB137").Value2
$row139 = $ws.Range("B139:AB139").Value2
$row141 = $ws.Range("B141:AB141").Value2
$row142 = $ws.Range("B142:AB142").Value2

# Apply cyclic rotation within each group: new[g[i]] = old[g[i+1]] (wrap-around)
$ws.Range("B41:AB41").Value2 = $row42
$ws.Range("B42:AB42").Value2 = $row43
$ws.Range("B43:AB43").Value2 = $row41

$ws.Range("B53:AB53").Value2 = $row54
$ws.Range("B54:AB54").Value2 = $row53

$ws.Range("B80:AB80").Value2 = $row81
$ws.Range("B81:AB81").Value2 = $row80

$ws.Range("B93:AB93").Value2 = $row94
$ws.Range("B94:AB94").Value2 = $row93

$ws.Range("B98:AB98").Value2 = $row99
$ws.Range("B99:AB99").Value2 = $row98

$ws.Range("B102:AB102").Value2 = $row103
$ws.Range("B103:AB103").Value2 = $row102

$ws.Range("B111:AB111").Value2 = $row112
$ws.Range("B112:AB112").Value2 = $row111

$ws.Range("B123:AB123").Value2 = $row124
$ws.Range("B124:AB124").Value2 = $row123

$ws.Range("B130:AB130").Value2 = $row132
$ws.Range("B132:AB132").Value2 = $row130

$ws.Range("B137:AB137").Value2 = $row139
$ws.Range("B139:AB139").Value2 = $row137

$ws.Range("B141:AB141").Value2 = $row142
$ws.Range("B142:AB142").Value2 = $row141

